$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$c = $t.Cell(1, 1)
$c.Range.Text = "56 x 28" + $nl + "  2    8" + $nl + "  ----" + $nl + "5|    |" + $nl + "6|    |"

$c = $t.Cell(1, 2)
$c.Range.Text = "30 x 71" + $nl + "  7    1" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"

$c = $t.Cell(1, 3)
$c.Range.Text = "53 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"

$c = $t.Cell(2, 1)
$c.Range.Text = "15 x 99" + $nl + "  9    9" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"

$c = $t.Cell(2, 2)
$c.Range.Text = "25 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "5|    |"

$c = $t.Cell(2, 3)
$c.Range.Text = "39 x 54" + $nl + "  5    4" + $nl + "  ----" + $nl + "3|    |" + $nl + "9|    |"

$c = $t.Cell(3, 1)
$c.Range.Text = "63 x 31" + $nl + "  3    1" + $nl + "  ----" + $nl + "6|    |" + $nl + "3|    |"

$c = $t.Cell(3, 2)
$c.Range.Text = "18 x 30" + $nl + "  3    0" + $nl + "  ----" + $nl + "1|    |" + $nl + "8|    |"

$c = $t.Cell(3, 3)
$c.Range.Text = "78 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "8|    |"

$c = $t.Cell(4, 1)
$c.Range.Text = "25 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "2|    |" + $nl + "5|    |"

$c = $t.Cell(4, 2)
$c.Range.Text = "51 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "5|    |" + $nl + "1|    |"

$c = $t.Cell(4, 3)
$c.Range.Text = "21 x 42" + $nl + "  4    2" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"

$c = $t.Cell(5, 1)
$c.Range.Text = "45 x 55" + $nl + "  5    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "5|    |"

$c = $t.Cell(5, 2)
$c.Range.Text = "19 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"

$c = $t.Cell(5, 3)
$c.Range.Text = "29 x 79" + $nl + "  7    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
